# Add a new "21-jun" data column (M) to the worksheet, mirroring the
# existing daily-tracking columns (C..L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("M1").Value = "21-jun"

# New values for the data rows (row -> value)
$values = @{
    2  = 0
    3  = 13.17106684475379
    4  = 16.880024755672888
    5  = 16.111859843455196
    6  = 0
    7  = 9.3061338130800877
    8  = 7.3290334555807117
    9  = 16.140204137518531
    10 = 12.748526836025437
    11 = 12.781254427131508
    12 = 0
    13 = 12.717347754040736
    14 = 0
    15 = 0
    16 = 16.034844474149928
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 13).Value = $values[$row]
}

# Update the active selection to match the authored workbook
$ws.Range("O6").Select()
